# Insert a new weekly record at row 55 (Zanahoria, Vega Monumental Concepción),
# pushing the existing rows 55-134 down to 56-135.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows(55).Insert()

$ws.Cells.Item(55, 1).Value  = 11
$ws.Cells.Item(55, 2).Value  = "Vega Monumental Concepción"
$ws.Cells.Item(55, 3).Value  = "Bíobío"
$ws.Cells.Item(55, 4).Value  = 44482
$ws.Cells.Item(55, 5).Value  = 8
$ws.Cells.Item(55, 6).Value  = 100114013
$ws.Cells.Item(55, 7).Value  = "Zanahoria"
$ws.Cells.Item(55, 8).Value  = "Sin especificar"
$ws.Cells.Item(55, 9).Value  = "Primera"
$ws.Cells.Item(55, 10).Value = 250
$ws.Cells.Item(55, 11).Value = 7500
$ws.Cells.Item(55, 12).Value = 8000
$ws.Cells.Item(55, 13).Value = 7700
$ws.Cells.Item(55, 14).Value = "$/saco 20 kilos"
$ws.Cells.Item(55, 15).Value = "Chillán"
$ws.Cells.Item(55, 16).Value = 385
$ws.Cells.Item(55, 17).Value = 20
$ws.Cells.Item(55, 18).Value = "Hortaliza"
